$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("CN1").Value = "test"
